$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("account")

# --- Insert two new columns (birthdate, city) before the existing "email" column ---
$ws.Range("E1:F1").EntireColumn.Insert()

# --- Headers ---
$ws.Cells.Item(1,5).Value = "birthdate"
$ws.Cells.Item(1,6).Value = "city"

# --- Row 2 (Admin) --- set city before birthdate so shared-string order matches
$ws.Cells.Item(2,6).Value = "'Jambi"
$ws.Cells.Item(2,5).Value = "'2002-12-03"

# --- Row 3 (User) ---
$ws.Cells.Item(3,6).Value = "'Jambi"
$ws.Cells.Item(3,5).Value = "'2003-12-03"

# --- Column widths for the two new columns ---
$ws.Columns.Item(5).ColumnWidth = 16.92
$ws.Columns.Item(6).ColumnWidth = 13.67

# --- Fix up hyperlinks that used to sit on the (now-shifted) email column ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2,7), "mailto:admin@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3,7), "mailto:user@gmail.com") | Out-Null

# Adding the hyperlink re-applies a style; put the original "Hyperlink" cell style back
$ws.Cells.Item(2,7).Style = "Hyperlink"
$ws.Cells.Item(3,7).Style = "Hyperlink"

# --- Make "account" the selected/active tab (was "contentCategory") ---
$ws.Activate()
$ws.Range("E4").Select() | Out-Null
